$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of row -> column -> new value
$updates = @{
    13  = @{ J = -0.0165; K = -0.0037; L = 0.1672;  M = 0.0786;  N = 0.0371;  O = -0.0128; P = -0.0628; Q = 0.0278;  R = 0.0256 }
    15  = @{ J = -0.4514; K = -0.2684; L = -0.1513; M = -0.0996; N = -0.0022; O = 0.2975;  P = 0.3056;  Q = 0.3289;  R = 0.3225 }
    23  = @{ J = -0.5935; K = -0.1412; L = -0.015;  M = 0.0539;  N = -0.0446; O = -0.0601; P = -0.0578; Q = -0.0556; R = -0.0351 }
    31  = @{ J = -1.9323; K = -1.6849; L = -2.403;  M = -0.8625; N = -0.4285; O = -0.1597; P = -0.3649; Q = -0.3748; R = -0.0532 }
    47  = @{ J = -0.0159; K = 0.0745;  L = -0.0054; M = 0.1739;  N = 0.1421;  O = 0.1185;  P = 0.0954;  Q = -0.0667; R = -0.0538 }
    69  = @{ J = 0.0071;  K = 0.0148;  L = 0.0227;  M = 0.0295;  N = 0.0424;  O = 0.0357;  P = 0.0288;  Q = 0.0187;  R = 0.0048 }
    71  = @{ J = 0.236;   K = 0.1518;  L = 0.1694;  M = 0.1514;  N = 0.1347;  O = 0.1275;  P = 0.1296;  Q = 0.0063;  R = 0.0168 }
    79  = @{ J = 0.1852;  K = 0.1992;  L = 0.1946;  M = 0.0499;  N = 0.0092;  O = -0.0089; P = -0.0087; Q = -0.0085; R = -0.0083 }
    87  = @{ J = 0.3982;  K = 0.3253;  L = 0.3624;  M = 0.2401;  N = 0.2028;  O = 0.1745;  P = 0.1732;  Q = -0.0016; R = 0.0174 }
    103 = @{ J = -0.0391; K = -0.0414; L = -0.0439; M = -0.0147; N = -0.011;  O = -0.0073; P = -0.0036; Q = 0.0032;  R = 0.0029 }
}

foreach ($row in $updates.Keys) {
    $cols = $updates[$row]
    foreach ($col in $cols.Keys) {
        $ws.Range("$col$row").Value = $cols[$col]
    }
}
